$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 4217.35
$ws.Range("I6").Value = 4630.3887
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 13891.1661
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -13779.1661
$ws.Range("N6").Value = -1724
$ws.Range("H55").Value = 183.46666
$ws.Range("I55").Value = 118
$ws.Range("J55").Value = 314.4
$ws.Range("K55").Value = 118
$ws.Range("L55").Value = 314.4
$ws.Range("M55").Value = 96
$ws.Range("N55").Value = -742.4
$ws.Range("H92").Value = 784.8889
$ws.Range("I92").Value = 115.25
$ws.Range("J92").Value = 2124.1667
$ws.Range("K92").Value = 115.25
$ws.Range("L92").Value = 2124.1667
$ws.Range("M92").Value = 1132.75
$ws.Range("N92").Value = -4620.1667
$ws.Range("H96").Value = 438.3684
$ws.Range("I96").Value = 450.66666
$ws.Range("K96").Value = 1351.99998
$ws.Range("M96").Value = 21.00001999999995
$ws.Range("H97").Value = 230
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 275.25
$ws.Range("I99").Value = 275.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 825.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 672.25
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 1558
$ws.Range("I100").Value = 1495
$ws.Range("K100").Value = 1495
$ws.Range("M100").Value = -954
$ws.Range("H101").Value = 14286493
$ws.Range("I101").Value = 28571764
$ws.Range("J101").Value = 1221.4286
$ws.Range("K101").Value = 85715292
$ws.Range("L101").Value = 3664.2858
$ws.Range("M101").Value = -85713670
$ws.Range("N101").Value = -6908.2858
$ws.Range("H106").Value = 1634.2222
$ws.Range("I106").Value = 1463.6111
$ws.Range("J106").Value = 1975.4445
$ws.Range("K106").Value = 1463.6111
$ws.Range("L106").Value = 1975.4445
$ws.Range("M106").Value = -832.6111000000001
$ws.Range("N106").Value = -3237.4445
$ws.Range("H129").Value = 993.3469
$ws.Range("I129").Value = 496.66666
$ws.Range("J129").Value = 1025.7391
$ws.Range("K129").Value = 1489.99998
$ws.Range("L129").Value = 3077.2173
$ws.Range("M129").Value = 3510.00002
$ws.Range("N129").Value = -13077.2173

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50014900
$ws.Range("J3").Value = 50014900
$ws.Range("L3").Value = 50014900
$ws.Range("N3").Value = -50015130
$ws.Range("H63").Value = 5926.5
$ws.Range("I63").Value = 1952.5
$ws.Range("J63").Value = 9900.5
$ws.Range("K63").Value = 1952.5
$ws.Range("L63").Value = 9900.5
$ws.Range("M63").Value = -1266.5
$ws.Range("N63").Value = -11272.5
$ws.Range("H66").Value = 5926.5
$ws.Range("I66").Value = 1952.5
$ws.Range("J66").Value = 9900.5
$ws.Range("K66").Value = 9762.5
$ws.Range("L66").Value = 49502.5
$ws.Range("M66").Value = -6330.5
$ws.Range("N66").Value = -56366.5
$ws.Range("H74").Value = 4870.625
$ws.Range("I74").Value = 5763.4736
$ws.Range("K74").Value = 5763.4736
$ws.Range("M74").Value = -4889.4736
$ws.Range("H77").Value = 4870.625
$ws.Range("I77").Value = 5763.4736
$ws.Range("K77").Value = 28817.368
$ws.Range("M77").Value = -24449.368
$ws.Range("H88").Value = 2592.3845
$ws.Range("I88").Value = 2434
$ws.Range("J88").Value = 2948.75
$ws.Range("K88").Value = 2434
$ws.Range("L88").Value = 2948.75
$ws.Range("M88").Value = -2028
$ws.Range("N88").Value = -3760.75
$ws.Range("H91").Value = 2592.3845
$ws.Range("I91").Value = 2434
$ws.Range("J91").Value = 2948.75
$ws.Range("K91").Value = 2434
$ws.Range("L91").Value = 2948.75
$ws.Range("M91").Value = -1030
$ws.Range("N91").Value = -5756.75
$ws.Range("H97").Value = 452.6
$ws.Range("I97").Value = 417.05264
$ws.Range("J97").Value = 565.1667
$ws.Range("K97").Value = 417.05264
$ws.Range("L97").Value = 565.1667
$ws.Range("M97").Value = 78.94736
$ws.Range("N97").Value = -1557.1667
$ws.Range("H102").Value = 4077.75
$ws.Range("I102").Value = 2100
$ws.Range("J102").Value = 4737
$ws.Range("K102").Value = 2100
$ws.Range("L102").Value = 4737
$ws.Range("M102").Value = -478
$ws.Range("N102").Value = -7981
$ws.Range("H132").Value = 1868.638
$ws.Range("I132").Value = 1780.7812
$ws.Range("J132").Value = 1976.7693
$ws.Range("K132").Value = 5342.3436
$ws.Range("L132").Value = 5930.3079
$ws.Range("M132").Value = -2812.3436
$ws.Range("N132").Value = -10990.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3227.2727
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377
$ws.Range("H89").Value = 3227.2727
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884
$ws.Range("H94").Value = 922.1177
$ws.Range("I94").Value = 989.5
$ws.Range("J94").Value = 862.2222
$ws.Range("K94").Value = 989.5
$ws.Range("L94").Value = 862.2222
$ws.Range("M94").Value = -538.5
$ws.Range("N94").Value = -1764.2222
$ws.Range("H99").Value = 2517.45
$ws.Range("I99").Value = 1455
$ws.Range("K99").Value = 1455
$ws.Range("M99").Value = 43
$ws.Range("H105").Value = 2953.57
$ws.Range("I105").Value = 1870
$ws.Range("J105").Value = 2987.0825
$ws.Range("K105").Value = 1870
$ws.Range("L105").Value = 2987.0825
$ws.Range("M105").Value = -123
$ws.Range("N105").Value = -6481.0825
$ws.Range("H134").Value = 1487.4615
$ws.Range("I134").Value = 1539.7273
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 4619.1819
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -2084.1819
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1220.9231
$ws.Range("I31").Value = 1457.2778
$ws.Range("J31").Value = 1095.7941
$ws.Range("K31").Value = 1457.2778
$ws.Range("L31").Value = 1095.7941
$ws.Range("M31").Value = -1162.2778
$ws.Range("N31").Value = -1685.7941
$ws.Range("H34").Value = 1220.9231
$ws.Range("I34").Value = 1457.2778
$ws.Range("J34").Value = 1095.7941
$ws.Range("K34").Value = 1457.2778
$ws.Range("L34").Value = 1095.7941
$ws.Range("M34").Value = -1255.2778
$ws.Range("N34").Value = -1499.7941
$ws.Range("H105").Value = 1002150.1
$ws.Range("I105").Value = 1668383.1
$ws.Range("J105").Value = 2800.5
$ws.Range("K105").Value = 1668383.1
$ws.Range("L105").Value = 2800.5
$ws.Range("M105").Value = -1666636.1
$ws.Range("N105").Value = -6294.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 747.7646999999999
$ws.Range("I4").Value = 367.8889
$ws.Range("J4").Value = 1175.125
$ws.Range("K4").Value = 1103.6667
$ws.Range("L4").Value = 3525.375
$ws.Range("M4").Value = -991.6667
$ws.Range("N4").Value = -3749.375
$ws.Range("H68").Value = 1353.2559
$ws.Range("J68").Value = 1318.1852
$ws.Range("L68").Value = 3954.5556
$ws.Range("N68").Value = -5576.5556
$ws.Range("H71").Value = 1353.2559
$ws.Range("J71").Value = 1318.1852
$ws.Range("L71").Value = 11863.6668
$ws.Range("N71").Value = -19975.6668
$ws.Range("H131").Value = 867.6316
$ws.Range("J131").Value = 1071.0714
$ws.Range("L131").Value = 3213.2142
$ws.Range("N131").Value = -13293.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18998824
$ws.Range("I11").Value = 10196455
$ws.Range("J11").Value = 51693332
$ws.Range("K11").Value = 10196455
$ws.Range("L11").Value = 51693332
$ws.Range("M11").Value = -10196316
$ws.Range("N11").Value = -51693610
$ws.Range("H97").Value = 770679.25
$ws.Range("I97").Value = 910630
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 910630
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -910134
$ws.Range("N97").Value = -1942

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 36666.668
$ws.Range("I93").Value = 36666.668
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 36666.668
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -35418.668
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 55557948
$ws.Range("I100").Value = 1320
$ws.Range("K100").Value = 1320
$ws.Range("M100").Value = -779

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1076.2
$ws.Range("I100").Value = 493.66666
$ws.Range("J100").Value = 1950
$ws.Range("K100").Value = 987.33332
$ws.Range("L100").Value = 3900
$ws.Range("M100").Value = -446.33332
$ws.Range("N100").Value = -4982

